$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.373.69"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").Value = "3.181.15"
$ws.Range("E3").Value = "  -7.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "3.179.21"
$ws.Range("E9").Value = "  -7.70%  "
$ws.Range("E10").Value = "  -5.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.397"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.34%  "
$ws.Range("D13").Value = "3.727.71"
$ws.Range("E13").Value = "  -7.81%  "
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "64.351.86"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D18").Value = "3.177.70"
$ws.Range("E18").Value = "  -7.68%  "
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.504"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.21%  "
$ws.Range("E26").Value = "  -3.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.14%  "
$ws.Range("E36").Value = "  -5.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "155.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.812"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.95%  "
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.28%  "
$ws.Range("D42").Value = "2.622.08"
$ws.Range("E42").Value = "  -5.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "332.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("E46").Value = "  -4.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0272"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  -0.04%  "

Write-Host "Updated cryptos list"
